# SOR Testing_PES NA Strategy & Mktg.xlsx - "Updated IPS AIP hipo turnover"
#
# This script:
#  1. Fixes up "Internal Fill Rate" (E/G.. row5/row6) values and adds a new
#     "Commit/Forecast" row (row 7) for the "Internal Fill Rate" block on
#     several location sheets.
#  2. Updates several "Professional Voluntary Turnover" ytd/month values.
#  3. Adds a brand-new location sheet "Piedras Negras Fasco Mexico" with the
#     same 23-column layout as the other sheets.

$wb = $excel.ActiveWorkbook

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")

function Set-PctRow {
    param($ws, [int]$row, $values)
    # $values is a hashtable keyed by column letter (E..W); only keys present
    # are written (others are left untouched). Percent number format is
    # applied to every written numeric cell so the displayed output matches
    # the existing "0.0%" columns.
    foreach ($col in $values.Keys) {
        $cell = $ws.Range($col + $row)
        $cell.Value = $values[$col]
        $cell.NumberFormat = "0.0%"
    }
}

# ---------------------------------------------------------------------------
# 1. Fort Wayne Indiana (sheet1) - Internal Fill Rate block
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Fort Wayne Indiana")

$ws1.Range("E5").Value = 1
$ws1.Range("E5").NumberFormat = "0.0%"
$ws1.Range("E6").Value = 1
$ws1.Range("E6").NumberFormat = "0.0%"

# New row 7: Commit/Forecast, Internal Fill Rate
$ws1.Range("A7").Value = "PES"
$ws1.Range("B7").Value = "PES NA Strategy & Mktg"
$ws1.Range("C7").Value = "Fort Wayne Indiana"
$ws1.Range("D7").Value = "Internal Fill Rate"
$ws1.Range("F7").Value = "Commit/Forecast"

Set-PctRow $ws1 7 @{
    "E" = 1
    "M" = 1; "N" = 1; "O" = 1; "P" = 1; "Q" = 1; "R" = 1
    "S" = 1; "T" = 1; "U" = 1; "V" = 1; "W" = 1
}
# Leave G7:L7 empty but percent-formatted (matches the other sheets' layout)
foreach ($col in @("G","H","I","J","K","L")) {
    $ws1.Range($col + "7").NumberFormat = "0.0%"
}

# ---------------------------------------------------------------------------
# 2. Grafton Wisconsin (sheet2) - Professional Voluntary Turnover
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Grafton Wisconsin")

Set-PctRow $ws2 2 @{ "E" = 0.0481 }
Set-PctRow $ws2 3 @{ "E" = 0.0481 }
Set-PctRow $ws2 4 @{
    "E" = 0.0481
    "G" = 0.0556
    "J" = 0.0518
    "M" = 0
    "N" = 0
    "O" = 0.00801666666666667
    "P" = 0.00801666666666667
    "Q" = 0.00801666666666667
    "R" = 0.02405
    "S" = 0.00801666666666667
    "T" = 0.00801666666666667
    "U" = 0.00801666666666667
    "V" = 0.02405
    "W" = 0.0962
}

# ---------------------------------------------------------------------------
# 3. Monterrey Rbm Mexico (sheet6) - Professional Voluntary Turnover
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Monterrey Rbm Mexico")

Set-PctRow $ws6 2 @{ "E" = 0.1087 }
Set-PctRow $ws6 3 @{ "E" = 0.1087 }
Set-PctRow $ws6 4 @{
    "E" = 0.1087
    "M" = 0.1
    "N" = 0.1031
    "O" = 0.0181166666666667
    "P" = 0.0181166666666667
    "Q" = 0.0181166666666667
    "R" = 0.05435
    "S" = 0.0181166666666667
    "T" = 0.0181166666666667
    "U" = 0.0181166666666667
    "V" = 0.05435
    "W" = 0.2174
}

# ---------------------------------------------------------------------------
# 4. Tipp City Ohio (sheet8) - Professional Voluntary Turnover
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Tipp City Ohio")

Set-PctRow $ws8 2 @{ "E" = 0.0952 }
Set-PctRow $ws8 3 @{ "E" = 0.0952 }
Set-PctRow $ws8 4 @{
    "E" = 0.0952
    "K" = 0.1
    "M" = 0
    "N" = 0.1
    "O" = 0.0158666666666667
    "P" = 0.0158666666666667
    "Q" = 0.0158666666666667
    "R" = 0.0476
    "S" = 0.0158666666666667
    "T" = 0.0158666666666667
    "U" = 0.0158666666666667
    "V" = 0.0476
    "W" = 0.1904
}

# ---------------------------------------------------------------------------
# 5. Milwaukee Pmc Hq Wisconsin (sheet10) - Internal Fill Rate + Turnover
# ---------------------------------------------------------------------------
$ws10 = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")

Set-PctRow $ws10 2 @{ "E" = 0.4 }
Set-PctRow $ws10 3 @{ "E" = 0.4 }
Set-PctRow $ws10 4 @{
    "E" = 0.4
    "K" = 0.5
    "M" = 0
    "N" = 0.4348
    "O" = 0.0666666666666667
    "P" = 0.0666666666666667
    "Q" = 0.0666666666666667
    "R" = 0.2
    "S" = 0.0666666666666667
    "T" = 0.0666666666666667
    "U" = 0.0666666666666667
    "V" = 0.2
    "W" = 0.8
}

Set-PctRow $ws10 5 @{ "E" = 0.5 }
Set-PctRow $ws10 6 @{ "E" = 0.5 }

# New row 7: Commit/Forecast, Internal Fill Rate
$ws10.Range("A7").Value = "PES"
$ws10.Range("B7").Value = "PES NA Strategy & Mktg"
$ws10.Range("C7").Value = "Milwaukee Pmc Hq Wisconsin"
$ws10.Range("D7").Value = "Internal Fill Rate"
$ws10.Range("F7").Value = "Commit/Forecast"

Set-PctRow $ws10 7 @{
    "E" = 0.5
    "H" = 1
    "J" = 1
    "M" = 0
    "N" = 0
    "O" = 0.5; "P" = 0.5; "Q" = 0.5; "R" = 0.5
    "S" = 0.5; "T" = 0.5; "U" = 0.5; "V" = 0.5; "W" = 0.5
}
foreach ($col in @("G","I","K","L")) {
    $ws10.Range($col + "7").NumberFormat = "0.0%"
}

# ---------------------------------------------------------------------------
# 6. New sheet: Piedras Negras Fasco Mexico
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add($null, $lastSheet)
$wsNew.Name = "Piedras Negras Fasco Mexico"

$headers = @("segment_function","division_function","location","cvd","ytd","data_source","Jan","Feb","Mar","Q1","Apr","May","Jun","Q2","Jul","Aug","Sep","Q3","Oct","Nov","Dec","Q4","FY")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsNew.Range($cols[$i] + "1").Value = $headers[$i]
}

$wsNew.Range("A2").Value = "PES"
$wsNew.Range("B2").Value = "PES NA Strategy & Mktg"
$wsNew.Range("C2").Value = "Piedras Negras Fasco Mexico"
$wsNew.Range("D2").Value = "Professional Voluntary Turnover"
$wsNew.Range("F2").Value = "Commit/Forecast"

Set-PctRow $wsNew 2 @{
    "E" = 0
    "G" = 0; "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0
    "M" = 0; "N" = 0; "O" = 0; "P" = 0; "Q" = 0; "R" = 0
    "S" = 0; "T" = 0; "U" = 0; "V" = 0; "W" = 0
}

# Re-select the first sheet so the workbook opens where it started.
$wb.Worksheets.Item("Fort Wayne Indiana").Activate()
